$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled header "fertlity_rate" -> "fertility_rate"
$ws.Range("C1").Value = "fertility_rate"

# Update the selected cell (view state) to match the author's final selection
$ws.Range("C15").Select()

# Minor column width tweaks captured by the author while reviewing the sheet
# (values chosen so the engine's stored width rounds to the closest match
# of the target widths 11.33203125 / 14.109375)
$ws.Columns("C").ColumnWidth = 10.5
$ws.Columns("E").ColumnWidth = 13.3
